$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: insert a new row at row 2 for "Albert Park" (pushes Bairnsdale.. down) ---
$ws.Rows.Item(2).Insert()
$ws.Range("A2:D2").ClearFormats()
$ws.Range("A2").Value = "Albert Park"
$ws.Range("B2").Value = "The Guilty Moose Cafe  143 Victoria Avenue, Albert Park VIC 3206"
$ws.Range("C2").Value = "30/12/20 1pm-1:30pm"
$ws.Range("D2").Value = "Case ate at cafe"

# --- Step 2: insert two new rows at rows 5-6 for "Camberwell" entries (before Caulfield) ---
$ws.Range("A5:A6").EntireRow.Insert()
$ws.Range("A5:D6").ClearFormats()
$ws.Range("A5").Value = "Camberwell"
$ws.Range("B5").Value = "Crown Nails  766 Riversdale Road, Camberwell VIC 3124"
$ws.Range("C5").Value = "30/12/20 1:30pm-2:30pm"
$ws.Range("D5").Value = "Case attended"

$ws.Range("A6").Value = "Camberwell"
$ws.Range("B6").Value = "Tao Dumplings  1 Evans Place, Camberwell VIC 3124"
$ws.Range("C6").Value = "30/12/20 12:30pm-1:30pm"
$ws.Range("D6").Value = "Case ate at restaurant"

# --- Step 3: remove the "Moorabbin" row (now at row 21) ---
$ws.Rows.Item(21).Delete()

# --- Step 4: normalise the remaining exposure-period dates from 4-digit to 2-digit years ---
$ws.Range("C4").Value = "24/12/20 10:00am-11:05am"
$ws.Range("C9").Value = "29/12/20 09:30am-10:45am"
$ws.Range("C10").Value = "29/12/20 7:00pm-9:30pm"
$ws.Range("C11").Value = "29/12/20 05:30pm-05:50pm"
$ws.Range("C12").Value = "29/12/20 03:30pm-04:30pm"
$ws.Range("C13").Value = "30/12/20 9:30am-10:30am"
$ws.Range("C14").Value = "29/12/20 2:30pm-5:50pm"
$ws.Range("C15").Value = "30/12/20 11:55am-12:30pm"
$ws.Range("C16").Value = "29/12/20 1:00pm-2:00pm"
$ws.Range("C17").Value = "26/12/20 5:30pm-7:30pm"
$ws.Range("C18").Value = "27/12/20 4:00pm-7:30pm"
$ws.Range("C19").Value = "28/12/20 8:00pm-9:30pm"
$ws.Range("C20").Value = "29/12/20 6:00pm-7:00pm"
$ws.Range("C22").Value = "25/12/20 12:00pm-02:30pm"
$ws.Range("C23").Value = "23/12/20 8:00pm-10:00pm"
